# Add team record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: copy formatting from the existing last header cell (AC1,
# style index 1: bold, centered, bordered) then overwrite with the new labels
# so the new headers match the look of the rest of the header row.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team's 2021 record (92-70-0) for every player row (2-61).
$ws.Range("AD2:AD61").Value = 92
$ws.Range("AE2:AE61").Value = 70
$ws.Range("AF2:AF61").Value = 0
